$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" (E column) values so the list appears in
# reverse order (2110 down to 2104), and swap the "Valor Mora" (F column)
# values that travel with the first/last row accordingly.
$ws.Range("E16").Value = "2110"
$ws.Range("E17").Value = "2109"
$ws.Range("E18").Value = "2108"
$ws.Range("E19").Value = "2107"
$ws.Range("E20").Value = "2106"
$ws.Range("E21").Value = "2105"
$ws.Range("E22").Value = "2104"

$ws.Range("F16").Value = 29260
$ws.Range("F17").Value = 36341
$ws.Range("F18").Value = 36341
$ws.Range("F19").Value = 36341
$ws.Range("F20").Value = 36341
$ws.Range("F21").Value = 36341
$ws.Range("F22").Value = 36341
